$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-07 09:46:31", 0.0012),
    @("2023-12-07 09:48:06", 0.005000000000000001),
    @("2023-12-07 09:49:23", 0.004600000000000001)
)

$startRow = 52
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
